$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 0) Remove the original "_GoBack" bookmark (currently sitting in paragraph 3,
#    between "...داشته " and "است.") before we add a differently-located
#    bookmark of the same name later on - two same-named bookmarks would make
#    later name lookups ambiguous.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# 1) Paragraph 3: merge "...داشته " and "است." into a single run.
# ---------------------------------------------------------------------------
$p3Start = $d.Paragraphs(3).Range.Start
$full3 = $d.Content.Text
$idx3 = $full3.IndexOf("داشته ", $p3Start)
$idx3EndPhrase = $full3.IndexOf("است.", $idx3)
$idx3End = $idx3EndPhrase + "است.".Length

$r3 = $d.Range($idx3, $idx3End)
# Route the replacement through a placeholder first: the destination text is
# character-for-character identical to the concatenation of the existing
# runs, so a direct assignment would be a no-op for the engine's run-diff
# and the underlying runs would not be merged.
$r3.Text = "TEMP_PLACEHOLDER_MERGE_2"
$full4 = $d.Content.Text
$idx4 = $full4.IndexOf("TEMP_PLACEHOLDER_MERGE_2")
$r4 = $d.Range($idx4, $idx4 + "TEMP_PLACEHOLDER_MERGE_2".Length)
$r4.Text = "داشته است."

# ---------------------------------------------------------------------------
# 2) Paragraph 1: collapse the run-per-phrase breakdown of
#    "سیج را" / "نباید" (underline) / " در قالب یک " / "سازمان" (bold) /
#    " نگریست بلکه بسیج یک " / "فرهنگ" (bold) / " است. ... رقابت داشت"
#    into a single plain run with the concatenated text.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$idxStart = $full.IndexOf("سیج را")
$idxEndPhrase = $full.IndexOf("رقابت داشت")
$idxEnd = $idxEndPhrase + "رقابت داشت".Length
$mergedText = "سیج را نباید در قالب یک سازمان نگریست بلکه بسیج یک فرهنگ است. بسیج، یادگاری است از امام امت (ره) در دورانی که ابرقدرت های غرب و شرق با هم رقابت داشت"

$r = $d.Range($idxStart, $idxEnd)
$r.Text = "TEMP_PLACEHOLDER_MERGE_1"
$full2 = $d.Content.Text
$idx2 = $full2.IndexOf("TEMP_PLACEHOLDER_MERGE_1")
$r2 = $d.Range($idx2, $idx2 + "TEMP_PLACEHOLDER_MERGE_1".Length)
$r2.Text = $mergedText

# ---------------------------------------------------------------------------
# 3) Paragraph 1: add a collapsed "_GoBack" bookmark right after the final
#    run (after the trailing ".") and before the paragraph mark.
#    A bookmark collapsed exactly at "end-of-paragraph-content" needs a
#    scratch character appended first (inserting one at the boundary makes
#    the anchor sit mid-paragraph instead of right on the paragraph-mark
#    edge), then the scratch character is removed again afterwards.
# ---------------------------------------------------------------------------
$p1EndPos = $d.Paragraphs(1).Range.End - 1
$scratch = $d.Range($p1EndPos, $p1EndPos)
$scratch.InsertAfter("Z")

$bmRange = $d.Range($p1EndPos, $p1EndPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$zRange = $d.Range($p1EndPos, $p1EndPos + 1)
$zRange.Text = ""

# ---------------------------------------------------------------------------
# 4) Section page margins: 1440 (1") -> 720 (0.5") twips on all four sides.
#    PageSetup margins are expressed in points (1 pt = 20 twips), so
#    720 twips == 36 pt.
# ---------------------------------------------------------------------------
$d.PageSetup.TopMargin = 36
$d.PageSetup.RightMargin = 36
$d.PageSetup.BottomMargin = 36
$d.PageSetup.LeftMargin = 36

Write-Output "edit complete"
